# Update the PSSM score matrix (B2:K21) with recomputed values
# from the supplemental-figures re-run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object "object[,]" 20,10

$data[0,0] = -18.29470599718708
$data[0,1] = 1.813493001798357
$data[0,2] = -18.29470599718708
$data[0,3] = -18.29470599718708
$data[0,4] = -18.29470599718708
$data[0,5] = -18.29470599718708
$data[0,6] = -18.29470599718708
$data[0,7] = -18.29470599718708
$data[0,8] = -18.29470599718708
$data[0,9] = -18.29470599718708

$data[1,0] = -18.29470599718708
$data[1,1] = -18.29470599718708
$data[1,2] = -18.29470599718708
$data[1,3] = -18.29470599718708
$data[1,4] = -18.29470599718708
$data[1,5] = -18.29470599718708
$data[1,6] = -18.29470599718708
$data[1,7] = 2.344776917611177
$data[1,8] = -18.29470599718708
$data[1,9] = -18.29470599718708

$data[2,0] = -18.29470599718708
$data[2,1] = 2.098706457501396
$data[2,2] = 2.075518447418087
$data[2,3] = -18.29470599718708
$data[2,4] = 3.628572517510871
$data[2,5] = -18.29470599718708
$data[2,6] = 1.762012412189934
$data[2,7] = -18.29470599718708
$data[2,8] = 1.750724398787792
$data[2,9] = -18.29470599718708

$data[3,0] = -18.29470599718708
$data[3,1] = 1.890078564393683
$data[3,2] = -18.29470599718708
$data[3,3] = -18.29470599718708
$data[3,4] = -18.29470599718708
$data[3,5] = -18.29470599718708
$data[3,6] = -18.29470599718708
$data[3,7] = -18.29470599718708
$data[3,8] = -18.29470599718708
$data[3,9] = -18.29470599718708

$data[4,0] = -18.29470599718708
$data[4,1] = -18.29470599718708
$data[4,2] = -18.29470599718708
$data[4,3] = -18.29470599718708
$data[4,4] = -18.29470599718708
$data[4,5] = -18.29470599718708
$data[4,6] = -18.29470599718708
$data[4,7] = -18.29470599718708
$data[4,8] = -18.29470599718708
$data[4,9] = -18.29470599718708

$data[5,0] = 2.609761783229724
$data[5,1] = -18.29470599718708
$data[5,2] = -18.29470599718708
$data[5,3] = -18.29470599718708
$data[5,4] = -18.29470599718708
$data[5,5] = -18.29470599718708
$data[5,6] = -18.29470599718708
$data[5,7] = -18.29470599718708
$data[5,8] = -18.29470599718708
$data[5,9] = -18.29470599718708

$data[6,0] = -18.29470599718708
$data[6,1] = -18.29470599718708
$data[6,2] = -18.29470599718708
$data[6,3] = 1.9279957632031
$data[6,4] = -18.29470599718708
$data[6,5] = -18.29470599718708
$data[6,6] = -18.29470599718708
$data[6,7] = -18.29470599718708
$data[6,8] = -18.29470599718708
$data[6,9] = -18.29470599718708

$data[7,0] = 3.796588946091864
$data[7,1] = -18.29470599718708
$data[7,2] = -18.29470599718708
$data[7,3] = -18.29470599718708
$data[7,4] = -18.29470599718708
$data[7,5] = -18.29470599718708
$data[7,6] = -18.29470599718708
$data[7,7] = -18.29470599718708
$data[7,8] = -18.29470599718708
$data[7,9] = -18.29470599718708

$data[8,0] = -18.29470599718708
$data[8,1] = -18.29470599718708
$data[8,2] = -18.29470599718708
$data[8,3] = -18.29470599718708
$data[8,4] = -18.29470599718708
$data[8,5] = -18.29470599718708
$data[8,6] = -18.29470599718708
$data[8,7] = 0.9359923256305337
$data[8,8] = -18.29470599718708
$data[8,9] = 1.816400374121601

$data[9,0] = -18.29470599718708
$data[9,1] = -18.29470599718708
$data[9,2] = -18.29470599718708
$data[9,3] = 2.874989675828215
$data[9,4] = -18.29470599718708
$data[9,5] = 4.321923832592223
$data[9,6] = -18.29470599718708
$data[9,7] = -18.29470599718708
$data[9,8] = -18.29470599718708
$data[9,9] = 2.274983687224003

$data[10,0] = -18.29470599718708
$data[10,1] = -18.29470599718708
$data[10,2] = -18.29470599718708
$data[10,3] = -18.29470599718708
$data[10,4] = -18.29470599718708
$data[10,5] = -18.29470599718708
$data[10,6] = -18.29470599718708
$data[10,7] = -18.29470599718708
$data[10,8] = -18.29470599718708
$data[10,9] = -18.29470599718708

$data[11,0] = -18.29470599718708
$data[11,1] = -18.29470599718708
$data[11,2] = -18.29470599718708
$data[11,3] = 2.483015658791246
$data[11,4] = -18.29470599718708
$data[11,5] = -18.29470599718708
$data[11,6] = -18.29470599718708
$data[11,7] = -18.29470599718708
$data[11,8] = 1.641185616466374
$data[11,9] = 1.849669694845572

$data[12,0] = -18.29470599718708
$data[12,1] = -18.29470599718708
$data[12,2] = 1.160607174369407
$data[12,3] = -18.29470599718708
$data[12,4] = -18.29470599718708
$data[12,5] = -18.29470599718708
$data[12,6] = -18.29470599718708
$data[12,7] = -18.29470599718708
$data[12,8] = -18.29470599718708
$data[12,9] = 2.15469356455535

$data[13,0] = -18.29470599718708
$data[13,1] = -18.29470599718708
$data[13,2] = 1.039412313703489
$data[13,3] = -18.29470599718708
$data[13,4] = -18.29470599718708
$data[13,5] = -18.29470599718708
$data[13,6] = -18.29470599718708
$data[13,7] = -18.29470599718708
$data[13,8] = -18.29470599718708
$data[13,9] = -18.29470599718708

$data[14,0] = -18.29470599718708
$data[14,1] = -18.29470599718708
$data[14,2] = -18.29470599718708
$data[14,3] = -18.29470599718708
$data[14,4] = -18.29470599718708
$data[14,5] = -18.29470599718708
$data[14,6] = -18.29470599718708
$data[14,7] = -18.29470599718708
$data[14,8] = 1.621460671602772
$data[14,9] = -18.29470599718708

$data[15,0] = -18.29470599718708
$data[15,1] = 2.282528266632903
$data[15,2] = 2.683549889766591
$data[15,3] = -18.29470599718708
$data[15,4] = -18.29470599718708
$data[15,5] = -18.29470599718708
$data[15,6] = 1.432531453213944
$data[15,7] = 2.666139461364508
$data[15,8] = 2.857935973515719
$data[15,9] = -18.29470599718708

$data[16,0] = -18.29470599718708
$data[16,1] = -18.29470599718708
$data[16,2] = -18.29470599718708
$data[16,3] = -18.29470599718708
$data[16,4] = -18.29470599718708
$data[16,5] = -18.29470599718708
$data[16,6] = 2.233092648621328
$data[16,7] = 1.532818136672594
$data[16,8] = 1.673022642529573
$data[16,9] = -18.29470599718708

$data[17,0] = -18.29470599718708
$data[17,1] = -18.29470599718708
$data[17,2] = 1.440322735867994
$data[17,3] = -18.29470599718708
$data[17,4] = -18.29470599718708
$data[17,5] = -18.29470599718708
$data[17,6] = 1.318610465781272
$data[17,7] = 1.074518580864233
$data[17,8] = -18.29470599718708
$data[17,9] = -18.29470599718708

$data[18,0] = -18.29470599718708
$data[18,1] = 0.7508857082355408
$data[18,2] = 1.236296538580238
$data[18,3] = -18.29470599718708
$data[18,4] = 2.931998409033084
$data[18,5] = -18.29470599718708
$data[18,6] = 1.685546565564311
$data[18,7] = 0.7312905380685633
$data[18,8] = -18.29470599718708
$data[18,9] = 1.840403293575802

$data[19,0] = -18.29470599718708
$data[19,1] = 0.9613298827233718
$data[19,2] = -18.29470599718708
$data[19,3] = 1.70843442519294
$data[19,4] = -18.29470599718708
$data[19,5] = -18.29470599718708
$data[19,6] = 1.806019214047321
$data[19,7] = -18.29470599718708
$data[19,8] = -18.29470599718708
$data[19,9] = -18.29470599718708

$ws.Range("B2:K21").Value = $data
Write-Host "Applied PSSM updates"